# Task #1199: Add auto loan calculation to sample data, including the data itself.
#
# 1. Remove the "Housing:Mortgage Interest" row from the sample data table - its
#    interest will now be computed from loan data attached to the mortgage
#    principal row instead of being tracked as its own standalone line item.
# 2. Add a new "Loan" column to Table1, and populate the mortgage-principal row
#    with the loan definition (interest category, principal amount, rate, term,
#    origination date) as a JSON blob.
# 3. Leave the selection positioned on the new column's first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleDataPattern")

# Row 3 is "Housing:Mortgage Interest" - delete the entire row so the table
# (and everything below) shifts up by one.
$ws.Rows.Item(3).Delete()

# Grow Table1 by one column to hold the new "Loan" field.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListColumns.Add() | Out-Null

# Name the new column and fill in the loan definition for the mortgage
# principal row (now row 2, right after the header row).
$ws.Range("J1").Value = "Loan"
$ws.Range("J2").Value = '{ "interest": "Housing:Mortgage Interest", "amount": 375000, "rate": 3, "term": 360, "origination": "1/1/2010" }'

# Match the author's final selection.
$ws.Activate()
$ws.Range("J3").Select()
